# Update "想去人数" (interest count) figures in column F across sheets,
# matching the scraped data refresh from commit 456a3b4.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 1109
$ws.Range("F3").Value  = 4791
$ws.Range("F5").Value  = 1936
$ws.Range("F6").Value  = 577
$ws.Range("F10").Value = 1181
$ws.Range("F13").Value = 526
$ws.Range("F14").Value = 2030
$ws.Range("F15").Value = 638
$ws.Range("F17").Value = 540
$ws.Range("F19").Value = 256
$ws.Range("F20").Value = 127
$ws.Range("F21").Value = 127
$ws.Range("F25").Value = 2556
$ws.Range("F29").Value = 1667
$ws.Range("F32").Value = 534
$ws.Range("F34").Value = 4412

# 演出 (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 376
$ws.Range("F5").Value = 2

# 本地生活 (Local Life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 1387
$ws.Range("F7").Value = 434

# 全部类型 (All Types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 1387
$ws.Range("F6").Value  = 434
$ws.Range("F7").Value  = 1109
$ws.Range("F8").Value  = 4791
$ws.Range("F9").Value  = 1936
$ws.Range("F10").Value = 577
$ws.Range("F11").Value = 376
$ws.Range("F15").Value = 1181
$ws.Range("F21").Value = 526
$ws.Range("F22").Value = 2030
$ws.Range("F23").Value = 638
$ws.Range("F25").Value = 540
$ws.Range("F27").Value = 256
$ws.Range("F29").Value = 127
$ws.Range("F30").Value = 127
$ws.Range("F38").Value = 2556
$ws.Range("F44").Value = 1667
$ws.Range("F46").Value = 534
$ws.Range("F49").Value = 4412

$wb.Save()
